$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVCRSbRIC")

# The old column K ("ISIC 20T21") is being split into two columns:
# new K = "ISIC 20", new L = "ISIC 21". Insert a new column before K so
# everything from K onward shifts right by one, then fix up the headers.
$ws.Range("K1").EntireColumn.Insert()

$ws.Range("K1").Value = "ISIC 20"
$ws.Range("L1").Value = "ISIC 21"
$ws.Range("K2").Value = 0
